$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shift existing D:K data to F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting from column F (previously column D) into new D:E columns
$ws.Range("F5:F102").Copy() | Out-Null
$ws.Range("D5:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate new D/E columns with the newly reported quarter data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 175200
$ws.Range("E8").Value = 93800
$ws.Range("D9").Value = 5300
$ws.Range("E9").Value = 4800
$ws.Range("D10").Value = 169900
$ws.Range("E10").Value = 89000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 145600
$ws.Range("E17").Value = 129700
$ws.Range("D18").Value = 29700
$ws.Range("E18").Value = -35900
$ws.Range("D20").Value = -9800
$ws.Range("E20").Value = -200
$ws.Range("D21").Value = 22500
$ws.Range("E21").Value = -33600
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 19900
$ws.Range("E23").Value = -36200
$ws.Range("D24").Value = 12900
$ws.Range("E24").Value = -900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 7000
$ws.Range("E26").Value = -35300
$ws.Range("D27").Value = -1000
$ws.Range("E27").Value = -14500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 9800
$ws.Range("E32").Value = 200
$ws.Range("D33").Value = -1000
$ws.Range("E33").Value = -14500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -1000
$ws.Range("E35").Value = -14500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 315800
$ws.Range("E41").Value = 195000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 103600
$ws.Range("E43").Value = 71900
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 11600
$ws.Range("E45").Value = 6200
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 561400
$ws.Range("E47").Value = 669100
$ws.Range("D48").Value = 43300
$ws.Range("E48").Value = 45100
$ws.Range("D49").Value = 22700
$ws.Range("E49").Value = 22700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 363100
$ws.Range("E52").Value = 367200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1447400
$ws.Range("E54").Value = 1421700
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("D58").Value = 62800
$ws.Range("E58").Value = 42400
$ws.Range("D59").Value = 496500
$ws.Range("E59").Value = 472900
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 290000
$ws.Range("E61").Value = 290100
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1876300
$ws.Range("E66").Value = 1852400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -3564700
$ws.Range("E72").Value = -3559300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -428900
$ws.Range("E76").Value = -430700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -1000
$ws.Range("E81").Value = -14500
$ws.Range("D83").Value = 2600
$ws.Range("E83").Value = 2500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 40400
$ws.Range("E89").Value = -60700
$ws.Range("D91").Value = -2200
$ws.Range("E91").Value = -1200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 86900
$ws.Range("E94").Value = -54800
$ws.Range("D96").Value = -3800
$ws.Range("E96").Value = -3800
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 1600
$ws.Range("E100").Value = 126500
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 128800
$ws.Range("E102").Value = 11100

# Row 58 anomaly: F:J become "NA" (no data reported) instead of shifted values
$ws.Range("F58").Value = "NA"
$ws.Range("G58").Value = "NA"
$ws.Range("H58").Value = "NA"
$ws.Range("I58").Value = "NA"
$ws.Range("J58").Value = "NA"
